$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - column F updates
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 20
$wsExpo.Range("F7").Value = 1695
$wsExpo.Range("F8").Value = 26
$wsExpo.Range("F11").Value = 1692
$wsExpo.Range("F13").Value = 88
$wsExpo.Range("F17").Value = 13
$wsExpo.Range("F21").Value = 479
$wsExpo.Range("F22").Value = 300
$wsExpo.Range("F24").Value = 231

# Sheet "全部类型" (sheet4) - column F updates
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 20
$wsAll.Range("F7").Value = 1695
$wsAll.Range("F9").Value = 26
$wsAll.Range("F12").Value = 1692
$wsAll.Range("F14").Value = 88
$wsAll.Range("F18").Value = 13
$wsAll.Range("F22").Value = 479
$wsAll.Range("F23").Value = 300
$wsAll.Range("F25").Value = 231
